# Added table to handle multiple BDOnTimes (as requested by Paul).
# Removed redundant fields/fields handled by other tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column header (row 1) rename: "Assembly Parts" -> "BlowdownSpecification On Times" ---
$ws.Range("W1").Value = "BlowdownSpecification On Times"

# --- Row header (row 23) rename: "Assembly Parts" -> "Blowdown Specification On Times" ---
$ws.Range("A23").Value = "Blowdown Specification On Times"

# --- Relationship value edits ---
# Equipment <-> Equipment self link now populated (was None)
$ws.Range("B2").Value = "1:N"
# Equipment <-> new table link cleared (was 1:N)
$ws.Range("W2").Value = "None"
# Blowdown Specifications <-> new table link populated (was None)
$ws.Range("W18").Value = "1:N"
# New table self link cleared (was 1:1)
$ws.Range("B23").Value = "None"
# New table <-> Blowdown Specifications link populated (was None)
$ws.Range("R23").Value = "1:1"

# --- Clear the now-unused trailing column Z on row 1 (used range shrinks to A1:Y23) ---
[void]$ws.Range("Z1").Clear()

# --- New column W width to fit the longer header text ---
$ws.Range("W1").EntireColumn.ColumnWidth = 29.3

# --- Re-highlight (green fill) the cells whose relationship is now populated ---
$greenCells = @("B2", "W18", "R23")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Interior.Color = 65280
}

# --- Remove highlight (clear fill) from the cells whose relationship is now "None" / cleared ---
$clearFillCells = @("W2", "H4", "B23", "H23")
foreach ($addr in $clearFillCells) {
    $ws.Range($addr).Interior.Pattern = -4142
}

# --- B23 loses its special numeric format, back to General ---
$ws.Range("B23").NumberFormat = "General"

# --- W7 picks up the same "Roboto / white fill" styling as its row neighbour V7 ---
$ws.Range("W7").Font.Name = "Roboto"
$ws.Range("W7").Font.Color = 255
$ws.Range("W7").Interior.Color = 16777215

# --- Move the active selection to W2, matching the saved view state ---
$ws.Range("W2").Select()

Write-Host "edit complete"
